$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (D, J, K, L, M, P) values per row, derived from the target diff.
$rows = @{
    2  = @(44497, 20, 4000, 4000, 4000, 4000)
    3  = @(44259, 30, 4000, 4000, 4000, 4000)
    4  = @(44508, 30, 4000, 4000, 4000, 4000)
    5  = @(44176, 10, 4000, 4000, 4000, 4000)
    6  = @(44509, 20, 4000, 4000, 4000, 4000)
    7  = @(44316, 20, 4000, 4000, 4000, 4000)
    8  = @(44313, 20, 4000, 4000, 4000, 4000)
    9  = @(44291, 35, 4000, 4000, 4000, 4000)
    10 = @(44315, 40, 4000, 4000, 4000, 4000)
    12 = @(44365, 55, 5000, 5000, 5000, 5000)
    13 = @(44504, 55, 4000, 4000, 4000, 4000)
    14 = @(44301, 40, 3000, 3000, 3000, 3000)
    15 = @(44312, 50, 4000, 4000, 4000, 4000)
    16 = @(44498, 40, 4000, 4000, 4000, 4000)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $vals[1]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals[2]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals[3]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals[4]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals[5]   # P - Precio $/Kg
}
